# Update column F (dSF) values for specific rows per repulled/pushed data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = -5
    5  = -3
    8  = 2
    14 = -5
    18 = 0
    19 = 2
    22 = -3
    27 = -7
    32 = -9
    33 = 0
    36 = 4
    37 = -4
    40 = 1
    47 = -4
    51 = 1
    59 = 1
    62 = -3
    67 = -3
    70 = -5
    74 = 1
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
